$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain text (their new values would otherwise
# be auto-converted to numbers by Excel, losing the exact text formatting).
$textCells = @("D4","D5","D6","D7","D8","D9","D11","D13","D14","D15","D16","D17","D19","D21","D25","D26","D28","D30","D31","D32","D33","D34","D35","D37","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50","D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values from the crypto price refresh
$ws.Range("D2").Value = "26.917.02"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.871.86"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "304.74"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.5096"
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("D8").Value = "0.3667"
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("D9").Value = "0.07183"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").Value = "20.66"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.878.76"
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.07493"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "94.49"
$ws.Range("E14").Value = "  +5.31%  "
$ws.Range("D15").Value = "5.227"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").Value = "0.9994"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "0.000008508"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").Value = "0.9994"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "26.961.41"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "5.015"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "2.115.03"
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").Value = "148.01"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").Value = "1.776"
$ws.Range("E26").Value = "  -3.42%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "2.084"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").Value = "4.700"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").Value = "4.717"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").Value = "0.09164"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "0.05056"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").Value = "0.7495"
$ws.Range("E34").Value = "  +3.56%  "
$ws.Range("D35").Value = "2.980"
$ws.Range("E35").Value = "  -2.97%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").Value = "3.218"
$ws.Range("E37").Value = "  +4.30%  "
$ws.Range("D38").Value = "2.524"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").Value = "0.5643"
$ws.Range("E39").Value = "  +6.95%  "
$ws.Range("D40").Value = "0.01990"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").Value = "1.071"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").Value = "6.617"
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("D43").Value = "115.60"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "8.548"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D46").Value = "0.4786"
$ws.Range("E46").Value = "  +3.71%  "
$ws.Range("D47").Value = "0.9991"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").Value = "10.10"
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("D49").Value = "1.559"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "36.98"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").Value = "63.11"
$ws.Range("E51").Value = "  -0.38%  "
